$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 12 (pushes the existing rows 12-20 down to 13-21).
# The new row inherits formatting (styles) from the row above, with no explicit height.
$ws.Rows.Item(12).Insert()

# Fill in the newly inserted row 12: only column A gets a label, B/C stay empty.
$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B12:C12").Clear()

# Row 13 (previously row 12, "Programa resumido:") - update the value shown in B/C.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (previously row 14, "Programa:") - B/C content replaced with this date-like
# text. Assigning the literal string directly would make Excel auto-convert it to a
# real date value, so instead put it in via a text formula and then paste the
# evaluated result back as a plain value (keeps it as text, keeps the original style).
$ws.Range("B15").Formula = '="01/01/2020"'
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").Formula = '="01/01/2020"'
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# Row 18 (previously row 17, "Método:") - B/C content replaced with this text.
$ws.Range("B18").Value = "5840820 - Gustavo Aristides Santana Martinez"
$ws.Range("C18").Value = "5840820 - Gustavo Aristides Santana Martinez"

# Row 10 ("Objetivos:") - the long Portuguese objectives paragraph is replaced.
$ws.Range("B10").Value = "5840820 - Gustavo Aristides Santana Martinez"
$ws.Range("C10").Value = "5840820 - Gustavo Aristides Santana Martinez"
